$wb = $excel.ActiveWorkbook

# Work on the third worksheet ("RO & CO Hearing Allocation")
$ws = $wb.Worksheets.Item(3)

# Rename the sheet
$ws.Name = "RO Allocations"

# Update the title text in A1
$ws.Range("A1").Value = "Allocation of Regional Office Video Hearings"

# Delete row 4 (the "Central Office" allocation row), shifting all
# subsequent rows up by one and removing the now-unused trailing row.
$ws.Rows("4").Delete()
